# Update column G ("K") values on Sheet1 to reflect recalculated strike counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 11
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 1
